# Insert a new data row at row 99 (pushing existing rows 99-141 down to 100-142)
# and populate it with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 99; Excel shifts rows 99:141 down to 100:142
# and copies the number-format of the row below (row 100, formerly row 99) onto the new row.
$ws.Rows("99:99").Insert()

# Populate the newly inserted row 99 with its values.
$ws.Cells.Item(99, 1).Value2  = 7
$ws.Cells.Item(99, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(99, 3).Value2  = "Ñuble"
$ws.Cells.Item(99, 4).Value2  = 44452
$ws.Cells.Item(99, 5).Value2  = 16
$ws.Cells.Item(99, 6).Value2  = 100112043
$ws.Cells.Item(99, 7).Value2  = "Pepino ensalada"
$ws.Cells.Item(99, 8).Value2  = "Sin especificar"
$ws.Cells.Item(99, 9).Value2  = "Primera"
$ws.Cells.Item(99, 10).Value2 = 300
$ws.Cells.Item(99, 11).Value2 = 16000
$ws.Cells.Item(99, 12).Value2 = 17000
$ws.Cells.Item(99, 13).Value2 = 16500
$ws.Cells.Item(99, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(99, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(99, 16).Value2 = 275
$ws.Cells.Item(99, 17).Value2 = 60
$ws.Cells.Item(99, 18).Value2 = "Hortaliza"
